{"js": "// Replace the song's opening line to match the new lyric.\nconst body = context.document.body;\nconst searchResults = body.search(\"Pushpa Raj Taggede le\u2026!\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Pushpa taggede le\u2026!\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\"Pushpa Raj Taggede le\u2026!\", $false, $false, $false, $false, $false, $true, 1, $false, \"Pushpa taggede le\u2026!\", 2)\n"}
